$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("GUEVARA IDROGO DENNIS PERCY", 109),
    @("TANTALEAN BUSTAMANTE ESTALIN YOEL", 107),
    @("INCIO SANCHEZ PAOLA KATHERINE", 106),
    @("PEREZ LINARES TATHIANA", 98),
    @("LINARES PEREZ YANASELY", 97),
    @("MONDRAGON HERNANDEZ WILMER JUNIOR", 96),
    @("MEDINA TAPIA ANA YULI", 95),
    @("HUAYHUA VALDIVIA LUZ EXMILDA", 95),
    @("CAMPOS PEREZ YOVERLY", 95),
    @("DELGADO VASQUEZ FLOR MAGALY", 93),
    @("CHAVEZ VILLANUEVA SILVIA JANETH", 89),
    @("LOZADA ROJAS LUZ ELENA", 89),
    @("SOTO LOZANO LUZDINA", 85)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
